$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2576
$ws.Range("J51").Value = 1975
$ws.Range("L51").Value = 1975
$ws.Range("N51").Value = -2943
$ws.Range("H100").Value = 1839.6842
$ws.Range("I100").Value = 1646.2858
$ws.Range("J100").Value = 2381.2
$ws.Range("K100").Value = 1646.2858
$ws.Range("L100").Value = 2381.2
$ws.Range("M100").Value = -1105.2858
$ws.Range("N100").Value = -3463.2
$ws.Range("H101").Value = 5200
$ws.Range("I101").Value = 400
$ws.Range("K101").Value = 1200
$ws.Range("M101").Value = 422
$ws.Range("H111").Value = 4337.364
$ws.Range("I111").Value = 3984.2307
$ws.Range("J111").Value = 4847.4443
$ws.Range("K111").Value = 11952.6921
$ws.Range("L111").Value = 14542.3329
$ws.Range("M111").Value = -8885.6921
$ws.Range("N111").Value = -20676.3329
$ws.Range("H121").Value = 2590
$ws.Range("I121").Value = 380
$ws.Range("K121").Value = 1140
$ws.Range("M121").Value = 607
$ws.Range("H123").Value = 29999.8
$ws.Range("J123").Value = 29999.8
$ws.Range("L123").Value = 29999.8
$ws.Range("N123").Value = -39799.8
$ws.Range("H138").Value = 1996.4073
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 1996.4073
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 5989.2219
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -16269.2219

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 469591.47
$ws.Range("I32").Value = 493030.53
$ws.Range("J32").Value = 41828.5
$ws.Range("K32").Value = 493030.53
$ws.Range("L32").Value = 41828.5
$ws.Range("M32").Value = -492743.53
$ws.Range("N32").Value = -42402.5
$ws.Range("H86").Value = 71456430
$ws.Range("J86").Value = 71456430
$ws.Range("L86").Value = 71456430
$ws.Range("N86").Value = -71458802
$ws.Range("H89").Value = 71456430
$ws.Range("J89").Value = 71456430
$ws.Range("L89").Value = 214369290
$ws.Range("N89").Value = -214381146
$ws.Range("H128").Value = 33939.5
$ws.Range("J128").Value = 33939.5
$ws.Range("L128").Value = 33939.5
$ws.Range("N128").Value = -43899.5
$ws.Range("H132").Value = 3267.6365
$ws.Range("I132").Value = 1794.56
$ws.Range("J132").Value = 5205.8945
$ws.Range("K132").Value = 5383.68
$ws.Range("L132").Value = 15617.6835
$ws.Range("M132").Value = -2853.68
$ws.Range("N132").Value = -20677.6835

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 11500
$ws.Range("J6").Value = 11500
$ws.Range("L6").Value = 11500
$ws.Range("N6").Value = -11726
$ws.Range("H105").Value = 11366373
$ws.Range("I105").Value = 12502710
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 12502710
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -12500963
$ws.Range("N105").Value = -6494
$ws.Range("H122").Value = 32890
$ws.Range("J122").Value = 32890
$ws.Range("L122").Value = 32890
$ws.Range("N122").Value = -42690
$ws.Range("H131").Value = 39824
$ws.Range("J131").Value = 39824
$ws.Range("L131").Value = 39824
$ws.Range("N131").Value = -49904

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1390
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1390
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1390
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1964
$ws.Range("H20").Value = 44866.91
$ws.Range("J20").Value = 44866.91
$ws.Range("L20").Value = 44866.91
$ws.Range("N20").Value = -45338.91
$ws.Range("H30").Value = 44866.91
$ws.Range("J30").Value = 44866.91
$ws.Range("L30").Value = 44866.91
$ws.Range("N30").Value = -45048.91
$ws.Range("H31").Value = 6564.522
$ws.Range("I31").Value = 1243.4667
$ws.Range("J31").Value = 16541.5
$ws.Range("K31").Value = 1243.4667
$ws.Range("L31").Value = 16541.5
$ws.Range("M31").Value = -948.4666999999999
$ws.Range("N31").Value = -17131.5
$ws.Range("H34").Value = 6564.522
$ws.Range("I34").Value = 1243.4667
$ws.Range("J34").Value = 16541.5
$ws.Range("K34").Value = 1243.4667
$ws.Range("L34").Value = 16541.5
$ws.Range("M34").Value = -1041.4667
$ws.Range("N34").Value = -16945.5
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H105").Value = 1821.4286
$ws.Range("I105").Value = 1800
$ws.Range("J105").Value = 1850
$ws.Range("K105").Value = 1800
$ws.Range("L105").Value = 1850
$ws.Range("M105").Value = -53
$ws.Range("N105").Value = -5344
$ws.Range("H113").Value = 1390
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1390
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1390
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5730
$ws.Range("H128").Value = 44866.91
$ws.Range("J128").Value = 44866.91
$ws.Range("L128").Value = 44866.91
$ws.Range("N128").Value = -54826.91

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2577
$ws.Range("I70").Value = 1506
$ws.Range("K70").Value = 4518
$ws.Range("M70").Value = -4203
$ws.Range("H73").Value = 2577
$ws.Range("I73").Value = 1506
$ws.Range("K73").Value = 4518
$ws.Range("M73").Value = -3426
$ws.Range("H113").Value = 912.6177
$ws.Range("I113").Value = 618.375
$ws.Range("K113").Value = 1855.125
$ws.Range("M113").Value = 314.875
$ws.Range("H133").Value = 13195.429
$ws.Range("I133").Value = 6448.3335
$ws.Range("K133").Value = 19345.0005
$ws.Range("M133").Value = -14285.0005
$ws.Range("H136").Value = 2549.875
$ws.Range("H139").Value = 4315.7188
$ws.Range("I139").Value = 1360
$ws.Range("K139").Value = 4080
$ws.Range("M139").Value = 1060
$ws.Range("H140").Value = 2102.1428
$ws.Range("J140").Value = 2493.75
$ws.Range("L140").Value = 7481.25
$ws.Range("N140").Value = -17841.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 29967590
$ws.Range("I80").Value = 56558836
$ws.Range("J80").Value = 52436.875
$ws.Range("K80").Value = 56558836
$ws.Range("L80").Value = 52436.875
$ws.Range("M80").Value = -56557838
$ws.Range("N80").Value = -54432.875
$ws.Range("H83").Value = 29967590
$ws.Range("I83").Value = 56558836
$ws.Range("J83").Value = 52436.875
$ws.Range("K83").Value = 282794180
$ws.Range("L83").Value = 262184.375
$ws.Range("M83").Value = -282789188
$ws.Range("N83").Value = -272168.375
$ws.Range("H113").Value = 112978.22
$ws.Range("I113").Value = 251447
$ws.Range("K113").Value = 251447
$ws.Range("M113").Value = -249277
$ws.Range("H114").Value = 31242.4
$ws.Range("J114").Value = 31242.4
$ws.Range("L114").Value = 31242.4
$ws.Range("N114").Value = -39920.4
$ws.Range("H122").Value = 1687.6111
$ws.Range("I122").Value = 1698.6471
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 5095.9413
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2645.9413
$ws.Range("N122").Value = -9400
$ws.Range("H124").Value = 47945
$ws.Range("J124").Value = 47945
$ws.Range("L124").Value = 47945
$ws.Range("N124").Value = -57765
$ws.Range("H136").Value = 21221
$ws.Range("J136").Value = 21221
$ws.Range("L136").Value = 63663
$ws.Range("N136").Value = -68763

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 67000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 67000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 67000
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -69246
$ws.Range("H90").Value = 67000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 67000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 201000
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -212232
$ws.Range("H132").Value = 2951.0557
$ws.Range("I132").Value = 1866.7778
$ws.Range("J132").Value = 4035.3333
$ws.Range("K132").Value = 5600.3334
$ws.Range("L132").Value = 12105.9999
$ws.Range("M132").Value = -3070.3334
$ws.Range("N132").Value = -17165.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 79325
$ws.Range("J86").Value = 79325
$ws.Range("L86").Value = 79325
$ws.Range("N86").Value = -81571
$ws.Range("H89").Value = 79325
$ws.Range("J89").Value = 79325
$ws.Range("L89").Value = 396625
$ws.Range("N89").Value = -407857
$ws.Range("H96").Value = 4647.469
$ws.Range("I96").Value = 3126
$ws.Range("J96").Value = 5087.8945
$ws.Range("K96").Value = 3126
$ws.Range("L96").Value = 5087.8945
$ws.Range("M96").Value = -1753
$ws.Range("N96").Value = -7833.8945
$ws.Range("H107").Value = 714.125
$ws.Range("I107").Value = 673.2857
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2019.8571
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -99.85710000000017
$ws.Range("N107").Value = -6840
$ws.Range("H123").Value = 22599.143
$ws.Range("J123").Value = 23638.8
$ws.Range("L123").Value = 23638.8
$ws.Range("N123").Value = -33438.8
